$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number/style formatting for new columns J, K, L from column I (same style per row)
for ($r = 1; $r -le 13; $r++) {
    $ws.Range("I$r").Copy() | Out-Null
    $ws.Range("J$r`:L$r").PasteSpecial(-4122) | Out-Null
}

# Set column widths for D..J to match target layout
$ws.Columns.Item(4).ColumnWidth = 8.142857142857142
$ws.Columns.Item(5).ColumnWidth = 10.714285714285714
$ws.Columns.Item(6).ColumnWidth = 11.0
$ws.Columns.Item(7).ColumnWidth = 8.142857142857142
$ws.Columns.Item(8).ColumnWidth = 13.142857142857142
$ws.Columns.Item(9).ColumnWidth = 9.714285714285714
$ws.Columns.Item(10).ColumnWidth = 8.142857142857142

# Update cell values
$ws.Range("D1").Value = "operator1"
$ws.Range("E1").Value = "rXpath1"
$ws.Range("F1").Value = "text1"
$ws.Range("G1").Value = "operator2"
$ws.Range("H1").Value = "rXpath2"
$ws.Range("I1").Value = "text2"
$ws.Range("J1").Value = "operator3"
$ws.Range("K1").Value = "subfunc3"
$ws.Range("L1").Value = "text3"
$ws.Range("D2").Value = "FILL"
$ws.Range("E2").Value = "rXpath1_web1"
$ws.Range("F2").Value = "text1_web1"
$ws.Range("G2").Value = "CLICK"
$ws.Range("H2").Value = "rXpath2_web1"
$ws.Range("I2").Value = "text2_web1"
$ws.Range("J2").Value = "FIND"
$ws.Range("K2").Value = "subfunc3_web1"
$ws.Range("L2").Value = "text3_web1"
$ws.Range("B3").Value = "yahoo"
$ws.Range("D3").Value = "FILL"
$ws.Range("E3").Value = "//input[@id='uh-search-box']"
$ws.Range("F3").Value = "what’s up? "
$ws.Range("G3").Value = "CLICK"
$ws.Range("H3").Value = "//button[@id='uh-search-button']"
$ws.Range("I3").Value = "text2_web2"
$ws.Range("J3").Value = "FIND"
$ws.Range("K3").Value = "YES"
$ws.Range("L3").Value = "what’s up? "
$ws.Range("D4").Value = "FILL"
$ws.Range("E4").Value = "//form[@name='frm_search']/input[@name='s']"
$ws.Range("F4").Value = "Christmas tree"
$ws.Range("G4").Value = "CLICK"
$ws.Range("H4").Value = "//form[@name='frm_search']/input[@value='search']"
$ws.Range("I4").Value = "text2_web3"
$ws.Range("J4").Value = "FIND"
$ws.Range("K4").Value = "NOT"
$ws.Range("L4").Value = "dkslfjksdh"
$ws.Range("B5").Value = "estate"
$ws.Range("C5").Value = "https://www.theestatesale.com/site/"
$ws.Range("D5").Value = "FILL"
$ws.Range("E5").Value = "//form[@name='frm_search']/input[@name='s']"
$ws.Range("F5").Value = "Christmas tree"
$ws.Range("G5").Value = "CLICK"
$ws.Range("H5").Value = "//form[@name='frm_search']/input[@value='search']"
$ws.Range("I5").Value = "text2_web3"
$ws.Range("J5").Value = "FIND"
$ws.Range("K5").Value = "YES"
$ws.Range("L5").Value = "Christmas tree"
$ws.Range("D6").Value = "FILL"
$ws.Range("E6").Value = "rXpath1_web5"
$ws.Range("F6").Value = "text1_web5"
$ws.Range("G6").Value = "CLICK"
$ws.Range("H6").Value = "rXpath2_web5"
$ws.Range("I6").Value = "text2_web5"
$ws.Range("J6").Value = "FIND"
$ws.Range("K6").Value = "subfunc3_text5"
$ws.Range("L6").Value = "text3_web5"
$ws.Range("D7").Value = "FILL"
$ws.Range("E7").Value = "rXpath1_web6"
$ws.Range("F7").Value = "text1_web6"
$ws.Range("G7").Value = "CLICK"
$ws.Range("H7").Value = "rXpath2_web6"
$ws.Range("I7").Value = "text2_web6"
$ws.Range("J7").Value = "FIND"
$ws.Range("K7").Value = "subfunc3_text6"
$ws.Range("L7").Value = "text3_web6"
$ws.Range("D8").Value = "FILL"
$ws.Range("E8").Value = "rXpath1_web7"
$ws.Range("F8").Value = "text1_web7"
$ws.Range("G8").Value = "CLICK"
$ws.Range("H8").Value = "rXpath2_web7"
$ws.Range("I8").Value = "text2_web7"
$ws.Range("J8").Value = "FIND"
$ws.Range("K8").Value = "subfunc3_text7"
$ws.Range("L8").Value = "text3_web7"
$ws.Range("D9").Value = "FILL"
$ws.Range("E9").Value = "rXpath1_web8"
$ws.Range("F9").Value = "text1_web8"
$ws.Range("G9").Value = "CLICK"
$ws.Range("H9").Value = "rXpath2_web8"
$ws.Range("I9").Value = "text2_web8"
$ws.Range("J9").Value = "FIND"
$ws.Range("K9").Value = "subfunc3_text8"
$ws.Range("L9").Value = "text3_web8"
$ws.Range("D10").Value = "FILL"
$ws.Range("E10").Value = "rXpath1_web9"
$ws.Range("F10").Value = "text1_web9"
$ws.Range("G10").Value = "CLICK"
$ws.Range("H10").Value = "rXpath2_web9"
$ws.Range("I10").Value = "text2_web9"
$ws.Range("J10").Value = "FIND"
$ws.Range("K10").Value = "subfunc3_text9"
$ws.Range("L10").Value = "text3_web9"
$ws.Range("D11").Value = "FILL"
$ws.Range("E11").Value = "rXpath1_web10"
$ws.Range("F11").Value = "text1_web10"
$ws.Range("G11").Value = "CLICK"
$ws.Range("H11").Value = "rXpath2_web10"
$ws.Range("I11").Value = "text2_web10"
$ws.Range("J11").Value = "FIND"
$ws.Range("K11").Value = "subfunc3_text10"
$ws.Range("L11").Value = "text3_web10"
$ws.Range("D12").Value = "FILL"
$ws.Range("E12").Value = "rXpath1_web11"
$ws.Range("F12").Value = "text1_web11"
$ws.Range("G12").Value = "CLICK"
$ws.Range("H12").Value = "rXpath2_web11"
$ws.Range("I12").Value = "text2_web11"
$ws.Range("J12").Value = "FIND"
$ws.Range("K12").Value = "subfunc3_text11"
$ws.Range("L12").Value = "text3_web11"
$ws.Range("D13").Value = "FILL"
$ws.Range("E13").Value = "rXpath1_web12"
$ws.Range("F13").Value = "text1_web12"
$ws.Range("G13").Value = "CLICK"
$ws.Range("H13").Value = "rXpath2_web12"
$ws.Range("I13").Value = "text2_web12"
$ws.Range("J13").Value = "FIND"
$ws.Range("K13").Value = "subfunc3_text12"
$ws.Range("L13").Value = "text3_web12"
